$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PPV_capacity")

$ws.Range("D39:D73").Value = 0
$ws.Range("D137:D169").Value = 0
$ws.Range("D232:D267").Value = 0
